$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: registration request
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "YCHERN"
$ws.Cells.Item(2, 3).Value = "ASFLI"
$ws.Cells.Item(2, 4).Value = "REGISTERPROJECT"
$ws.Cells.Item(2, 5).Value = "APPROVED"
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(2, 8).Style = "Normal"

# Row 3: deregistration request
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "YCHERN"
$ws.Cells.Item(3, 3).Value = "ASFLI"
$ws.Cells.Item(3, 4).Value = "DEREGISTERPROJECT"
$ws.Cells.Item(3, 5).Value = "APPROVED"
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(3, 8).Style = "Normal"

$wb.Save()
